# Insert a new "Note:" paragraph right after the existing paragraph that
# ends with "...Dr. Ioerger's heuristic." (the paragraph already containing
# the bold "Note: " run followed by the "Blue coloring implies..." text),
# and before the blank paragraph that precedes the results table.

$d = $word.ActiveDocument

# Locate the anchor paragraph by its distinctive trailing sentence.
$anchor = $d.Content
$found = $anchor.Find.Execute( `
    "Blue coloring implies that the plan length metric was better than Dr. Ioerger's heuristic. Red implies the plan length on the heuristic was worse than Dr. Ioerger's heuristic.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate anchor paragraph text."
}

# Collapse to the end of the found range and insert a brand-new, empty
# paragraph right after it (inherits the surrounding Times New Roman / 24pt
# paragraph-mark formatting, matching the rest of the document's notes).
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()

# Move into the freshly created paragraph.
$anchor.Collapse(0)
$anchor.Move(1, 1)

# Insert the full sentence as plain text first, remembering where it starts
# so we can go back and bold just the "Note:" label afterwards.
$noteStart = $anchor.Start
$anchor.InsertAfter("Note: After comparing results obtained from the server, it appears that the results obtained in the compute.cse server WILL DIFFER from the ones below, which were obtained from an M1 Macbook Air.")

# Bold only the "Note:" label (5 characters), leaving the rest of the
# sentence in regular weight, mirroring the other "Note:" paragraph above.
$labelRange = $d.Range($noteStart, $noteStart + 5)
$labelRange.Bold = 1

Write-Output "Inserted Note paragraph after results-color explanation."
